$wb = $excel.ActiveWorkbook

# ---- Copy "processed" -> "processed1" (appended at end of sheet list) ----
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$processedSrc = $wb.Worksheets.Item("processed")
$processedSrc.Copy($null, $last)
$processedNew = $wb.Worksheets.Item("processed (2)")
$processedNew.Name = "processed1"
# Update the listed sheet names so they point at the newly created sheets
$processedNew.Range("B1").Value = "ready_259781"
$processedNew.Range("B2").Value = "ready_285821"

# ---- Copy "ready_25978" -> "ready_259781" (appended at end of sheet list) ----
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ready1Src = $wb.Worksheets.Item("ready_25978")
$ready1Src.Copy($null, $last)
$ready1New = $wb.Worksheets.Item("ready_25978 (2)")
$ready1New.Name = "ready_259781"
$ready1New.PageSetup.LeftMargin = 50.4
$ready1New.PageSetup.RightMargin = 50.4
$ready1New.PageSetup.TopMargin = 54
$ready1New.PageSetup.BottomMargin = 54
$ready1New.PageSetup.HeaderMargin = 21.599999999999998
$ready1New.PageSetup.FooterMargin = 21.599999999999998

# ---- Copy "ready_28582" -> "ready_285821" (appended at end of sheet list) ----
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ready2Src = $wb.Worksheets.Item("ready_28582")
$ready2Src.Copy($null, $last)
$ready2New = $wb.Worksheets.Item("ready_28582 (2)")
$ready2New.Name = "ready_285821"
$ready2New.PageSetup.LeftMargin = 50.4
$ready2New.PageSetup.RightMargin = 50.4
$ready2New.PageSetup.TopMargin = 54
$ready2New.PageSetup.BottomMargin = 54
$ready2New.PageSetup.HeaderMargin = 21.599999999999998
$ready2New.PageSetup.FooterMargin = 21.599999999999998
